$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) cells remain plain text so values like
# "1.00" or "528.35" are not coerced into numbers, matching the
# original inline-string formatting used throughout this sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.027.69"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.062.39"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.35"
$ws.Range("E5").Value = "  +6.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.82"
$ws.Range("E6").Value = "  +6.65%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +5.36%  "
$ws.Range("E9").Value = "  +5.71%  "
$ws.Range("E10").Value = "  +7.24%  "
$ws.Range("E11").Value = "  +6.31%  "
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.588.88"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("E14").Value = "  +8.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  +16.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.983.33"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("E17").Value = "  +8.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.063.34"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.22"
$ws.Range("E19").Value = "  +7.42%  "
$ws.Range("E20").Value = "  +5.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.71"
$ws.Range("E21").Value = "  +5.13%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.70"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  +7.92%  "
$ws.Range("E25").Value = "  +5.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0978"
$ws.Range("E26").Value = "  +9.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.170"
$ws.Range("E27").Value = "  +4.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("E29").Value = "  +9.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.46"
$ws.Range("E30").Value = "  +10.21%  "
$ws.Range("E31").Value = "  +7.58%  "
$ws.Range("E32").Value = "  +6.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.27"
$ws.Range("E33").Value = "  +4.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").Value = "  +8.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.91"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.99"
$ws.Range("E36").Value = "  +7.26%  "
$ws.Range("E37").Value = "  +4.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.10"
$ws.Range("E38").Value = "  +12.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0701"
$ws.Range("E39").Value = "  +4.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.097.56"
$ws.Range("E40").Value = "  +2.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.85"
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("E42").Value = "  +11.47%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.667"
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.345.81"
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.05"
$ws.Range("E46").Value = "  +5.08%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.48"
$ws.Range("E47").Value = "  +6.20%  "
$ws.Range("E48").Value = "  +7.39%  "
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("E50").Value = "  +4.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.33"
$ws.Range("E51").Value = "  +6.92%  "
